$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Recorded By" column (G) lists the user(s)/systems that recorded a
# session. Each entry has been rotated so that the last name in the
# comma-separated list now appears first (i.e. a right-rotation of the list).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = @($val -split ", ")

        if ($parts.Count -gt 1) {
            $last = $parts[$parts.Count - 1]
            $rest = @($parts[0..($parts.Count - 2)])
            $newParts = @($last) + $rest
            $newVal = $newParts -join ", "

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
